$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1. "Team Name (optional)" row -> second cell gets "TEAM ONE"
$t.Cell(2, 2).Range.InsertAfter("TEAM ONE")

# 2. "Team members" cell: append "." to the end of each of the 4 member
#    paragraphs (Daniel Mutinda, Fredrick Odawa, Serah Wangari, Daniel Sifuna).
#    Locate them by scanning paragraphs between the "Team members" label
#    paragraph and the "Project Title" label paragraph - robust to any
#    earlier index shifts. NOTE: once a Table/Cell has been touched,
#    $d.Paragraphs.Item(n) becomes unreliable in this runtime, so we use
#    $d.Content.Paragraphs.Item(n) instead, which stays accurate.
$inTeamMembers = $false
$targets = @()
$count = $d.Content.Paragraphs.Count
for ($p = 1; $p -le $count; $p++) {
    $para = $d.Content.Paragraphs.Item($p)
    $txt = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Team members") {
        $inTeamMembers = $true
        continue
    }
    if ($txt -eq "Project Title") {
        $inTeamMembers = $false
    }
    if ($inTeamMembers -and $txt.Length -gt 0) {
        $targets += $p
    }
}

foreach ($idx in $targets) {
    $para = $d.Content.Paragraphs.Item($idx)
    $r = $para.Range
    # Collapse to just before the paragraph mark, then insert the period there
    $insertPoint = $d.Range($r.Start, $r.End - 1)
    $insertPoint.InsertAfter(".")
}

# 3. "Project Title" row -> second cell gets the project title text
$t2 = $d.Tables.Item(1)
$t2.Cell(5, 2).Range.InsertAfter("Home And Community Based Care(HCBC) HMIS")

# 4. "Short description (max. 25 words)" row -> second cell gets the
#    description text
$t3 = $d.Tables.Item(1)
$t3.Cell(6, 2).Range.InsertAfter("HCBC allows community health workers conduct data collection and surveys, provide referrals for free medication (ARV, TB) and treatment on behalf of the Ministry of Health.")
